$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1817.6364
$ws.Range("I2").Value = 1125
$ws.Range("K2").Value = 1125
$ws.Range("M2").Value = -1012
# Row 28
$ws.Range("H28").Value = 677.45
$ws.Range("I28").Value = 701.6111
$ws.Range("K28").Value = 701.6111
$ws.Range("M28").Value = -216.6111
# Row 80
$ws.Range("H80").Value = 45454790
$ws.Range("I80").Value = 50000212
$ws.Range("J80").Value = 600
$ws.Range("K80").Value = 150000636
$ws.Range("L80").Value = 1800
$ws.Range("M80").Value = -149999638
$ws.Range("N80").Value = -3796
# Row 83
$ws.Range("H83").Value = 45454790
$ws.Range("I83").Value = 50000212
$ws.Range("J83").Value = 600
$ws.Range("K83").Value = 450001908
$ws.Range("L83").Value = 5400
$ws.Range("M83").Value = -449996916
$ws.Range("N83").Value = -15384
# Row 86
$ws.Range("H86").Value = 22894.895
$ws.Range("I86").Value = 6499.857
$ws.Range("J86").Value = 32458.666
$ws.Range("K86").Value = 6499.857
$ws.Range("L86").Value = 32458.666
$ws.Range("M86").Value = -5376.857
$ws.Range("N86").Value = -34704.666
# Row 89
$ws.Range("H89").Value = 22894.895
$ws.Range("I89").Value = 6499.857
$ws.Range("J89").Value = 32458.666
$ws.Range("K89").Value = 32499.285
$ws.Range("L89").Value = 162293.33
$ws.Range("M89").Value = -26883.285
$ws.Range("N89").Value = -173525.33
# Row 112
$ws.Range("H112").Value = 2604.5
$ws.Range("J112").Value = 2805.818
$ws.Range("L112").Value = 8417.454000000002
$ws.Range("N112").Value = -10633.454
# Row 137
$ws.Range("H137").Value = 3018.875
$ws.Range("I137").Value = 2630.2
$ws.Range("K137").Value = 7890.599999999999
$ws.Range("M137").Value = -5340.599999999999
# Row 138
$ws.Range("H138").Value = 2578.875
$ws.Range("J138").Value = 5800
$ws.Range("L138").Value = 17400
$ws.Range("N138").Value = -27680
# Row 141
$ws.Range("H141").Value = 4205.3
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 19
$ws.Range("H19").Value = 12250
$ws.Range("I19").Value = 12250
$ws.Range("K19").Value = 12250
$ws.Range("M19").Value = -12021
# Row 37
$ws.Range("H37").Value = 24997.143
$ws.Range("J37").Value = 24997.143
$ws.Range("L37").Value = 24997.143
$ws.Range("N37").Value = -25543.143
# Row 97
$ws.Range("H97").Value = 2287.5
$ws.Range("I97").Value = 2287.5
$ws.Range("K97").Value = 2287.5
$ws.Range("M97").Value = -1791.5
# Row 110
$ws.Range("H110").Value = 4334.4443
$ws.Range("I110").Value = 1010
$ws.Range("K110").Value = 1010
$ws.Range("M110").Value = 1035
# Row 132
$ws.Range("H132").Value = 4473.3125
$ws.Range("I132").Value = 1720.2222
$ws.Range("J132").Value = 8013
$ws.Range("K132").Value = 5160.6666
$ws.Range("L132").Value = 24039
$ws.Range("M132").Value = -2630.6666
$ws.Range("N132").Value = -29099

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 1250374
$ws.Range("I7").Value = 1250374
$ws.Range("K7").Value = 1250374
$ws.Range("M7").Value = -1250261
# Row 82
$ws.Range("H82").Value = 21802.111
$ws.Range("I82").Value = 7243.8
$ws.Range("K82").Value = 7243.8
$ws.Range("M82").Value = -6860.8
# Row 85
$ws.Range("H85").Value = 21802.111
$ws.Range("I85").Value = 7243.8
$ws.Range("K85").Value = 7243.8
$ws.Range("M85").Value = -5917.8
# Row 105
$ws.Range("H105").Value = 2693.6
$ws.Range("I105").Value = 1992
$ws.Range("J105").Value = 5500
$ws.Range("K105").Value = 1992
$ws.Range("L105").Value = 5500
$ws.Range("M105").Value = -245
$ws.Range("N105").Value = -8994

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 69.166664
$ws.Range("I7").Value = 74.8
$ws.Range("K7").Value = 74.8
$ws.Range("M7").Value = 38.2
# Row 16
$ws.Range("H16").Value = 6000
$ws.Range("I16").Value = 10000
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -9713
$ws.Range("N16").Value = -2574
# Row 22
$ws.Range("H22").Value = 997
$ws.Range("I22").Value = 994.5
$ws.Range("J22").Value = 1002
$ws.Range("K22").Value = 994.5
$ws.Range("L22").Value = 1002
$ws.Range("M22").Value = -644.5
$ws.Range("N22").Value = -1702
# Row 31
$ws.Range("H31").Value = 4210
$ws.Range("I31").Value = 3449.6667
$ws.Range("K31").Value = 3449.6667
$ws.Range("M31").Value = -3154.6667
# Row 34
$ws.Range("H34").Value = 4210
$ws.Range("I34").Value = 3449.6667
$ws.Range("K34").Value = 3449.6667
$ws.Range("M34").Value = -3247.6667
# Row 41
$ws.Range("H41").Value = 17462.5
$ws.Range("J41").Value = 22000
$ws.Range("L41").Value = 22000
$ws.Range("N41").Value = -22856
# Row 50
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31250
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 59
$ws.Range("H59").Value = 28827.385
$ws.Range("J59").Value = 34999.832
$ws.Range("L59").Value = 34999.832
$ws.Range("N59").Value = -37289.832
# Row 60
$ws.Range("H60").Value = 18587.818
$ws.Range("J60").Value = 24997.857
$ws.Range("L60").Value = 24997.857
$ws.Range("N60").Value = -26019.857
# Row 99
$ws.Range("H99").Value = 4876.3335
$ws.Range("I99").Value = 6594.6665
$ws.Range("J99").Value = 3158
$ws.Range("K99").Value = 6594.6665
$ws.Range("L99").Value = 3158
$ws.Range("M99").Value = -5096.6665
$ws.Range("N99").Value = -6154
# Row 107
$ws.Range("H107").Value = 642.4375
$ws.Range("I107").Value = 575.1
$ws.Range("J107").Value = 754.6667
$ws.Range("K107").Value = 575.1
$ws.Range("L107").Value = 754.6667
$ws.Range("M107").Value = 1344.9
$ws.Range("N107").Value = -4594.6667
# Row 113
$ws.Range("H113").Value = 6000
$ws.Range("I113").Value = 10000
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -7830
$ws.Range("N113").Value = -6340
# Row 122
$ws.Range("H122").Value = 3140
$ws.Range("I122").Value = 3025.7646
$ws.Range("J122").Value = 3625.5
$ws.Range("K122").Value = 9077.293799999999
$ws.Range("L122").Value = 10876.5
$ws.Range("M122").Value = -6627.293799999999
$ws.Range("N122").Value = -15776.5
# Row 126
$ws.Range("H126").Value = 4876.3335
$ws.Range("I126").Value = 6594.6665
$ws.Range("J126").Value = 3158
$ws.Range("K126").Value = 19783.9995
$ws.Range("L126").Value = 9474
$ws.Range("M126").Value = -17313.9995
$ws.Range("N126").Value = -14414
# Row 132
$ws.Range("H132").Value = 1494.5
$ws.Range("I132").Value = 1494.5
$ws.Range("K132").Value = 4483.5
$ws.Range("M132").Value = -1953.5
# Row 134
$ws.Range("H134").Value = 15000
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3486.6924
$ws.Range("I122").Value = 3372.4443
$ws.Range("K122").Value = 10117.3329
$ws.Range("M122").Value = -7667.332900000001
# Row 136
$ws.Range("H136").Value = 3126.75
$ws.Range("I136").Value = 3126.75
$ws.Range("K136").Value = 9380.25
$ws.Range("M136").Value = -6830.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 876.5
$ws.Range("I6").Value = 730
$ws.Range("J6").Value = 949.75
$ws.Range("K6").Value = 730
$ws.Range("L6").Value = 949.75
$ws.Range("M6").Value = -615
$ws.Range("N6").Value = -1179.75
# Row 54
$ws.Range("H54").Value = 36472.125
$ws.Range("I54").Value = 8000
$ws.Range("K54").Value = 8000
$ws.Range("M54").Value = -7480
# Row 62
$ws.Range("H62").Value = 6000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 6000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# Row 122
$ws.Range("H122").Value = 3836.75
$ws.Range("I122").Value = 1449
$ws.Range("K122").Value = 4347
$ws.Range("M122").Value = -1897
# Row 132
$ws.Range("H132").Value = 1177.3334
$ws.Range("I132").Value = 1396
$ws.Range("J132").Value = 740
$ws.Range("K132").Value = 4188
$ws.Range("L132").Value = 2220
$ws.Range("M132").Value = -1658
$ws.Range("N132").Value = -7280
# Row 136
$ws.Range("H136").Value = 3077.8125
$ws.Range("I136").Value = 2949.6667
$ws.Range("K136").Value = 8849.000100000001
$ws.Range("M136").Value = -6299.000100000001

